$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 202.61539
$ws.Range("I9").Value = 205
$ws.Range("J9").Value = 189.5
$ws.Range("K9").Value = 205
$ws.Range("L9").Value = 189.5
$ws.Range("M9").Value = -36
$ws.Range("N9").Value = -527.5
$ws.Range("H33").Value = 126250470
$ws.Range("I33").Value = 1666800.9
$ws.Range("K33").Value = 1666800.9
$ws.Range("M33").Value = -1666571.9
$ws.Range("H34").Value = 4162.4287
$ws.Range("I34").Value = 4162.4287
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 4162.4287
$ws.Range("L34").Value = 0
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -3959.4287
$ws.Range("H36").Value = 4162.4287
$ws.Range("I36").Value = 4162.4287
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 4162.4287
$ws.Range("L36").Value = 0
$ws.Range("M36").ClearContents()
$ws.Range("N36").Value = -3447.4287
$ws.Range("H113").Value = 33060.5
$ws.Range("I113").Value = 48397
$ws.Range("K113").Value = 48397
$ws.Range("M113").Value = -45143
$ws.Range("H137").Value = 22729746
$ws.Range("I137").Value = 27780448
$ws.Range("J137").Value = 1587.25
$ws.Range("K137").Value = 83341344
$ws.Range("L137").Value = 4761.75
$ws.Range("M137").Value = -83338794
$ws.Range("N137").Value = -9861.75
$ws.Range("H141").Value = 5793.727
$ws.Range("I141").Value = 4966.5
$ws.Range("J141").Value = 7999.6665
$ws.Range("K141").Value = 14899.5
$ws.Range("L141").Value = 23998.9995
$ws.Range("M141").Value = -9719.5
$ws.Range("N141").Value = -34358.99950000001

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4299.6284
$ws.Range("I32").Value = 2484.5
$ws.Range("K32").Value = 2484.5
$ws.Range("M32").Value = -2197.5
$ws.Range("H61").Value = 3158.1428
$ws.Range("I61").Value = 2692.6365
$ws.Range("J61").Value = 4865
$ws.Range("K61").Value = 2692.6365
$ws.Range("L61").Value = 4865
$ws.Range("M61").Value = -2480.6365
$ws.Range("N61").Value = -5289
$ws.Range("H88").Value = 1938.0769
$ws.Range("I88").Value = 1127.1111
$ws.Range("K88").Value = 1127.1111
$ws.Range("M88").Value = -721.1111000000001
$ws.Range("H91").Value = 1938.0769
$ws.Range("I91").Value = 1127.1111
$ws.Range("K91").Value = 1127.1111
$ws.Range("M91").Value = 276.8888999999999
$ws.Range("H122").Value = 5059.7334
$ws.Range("I122").Value = 4511
$ws.Range("K122").Value = 13533
$ws.Range("M122").Value = -11083
$ws.Range("H132").Value = 7187.3213
$ws.Range("I132").Value = 3159.1667
$ws.Range("K132").Value = 9477.500100000001
$ws.Range("M132").Value = -6947.500100000001
$ws.Range("H136").Value = 3158.1428
$ws.Range("I136").Value = 2692.6365
$ws.Range("J136").Value = 4865
$ws.Range("K136").Value = 8077.9095
$ws.Range("L136").Value = 14595
$ws.Range("M136").Value = -5527.9095
$ws.Range("N136").Value = -19695

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 6805.8887
$ws.Range("I107").Value = 7037.5713
$ws.Range("K107").Value = 7037.5713
$ws.Range("M107").Value = -5117.5713
$ws.Range("H110").Value = 11950
$ws.Range("J110").Value = 11950
$ws.Range("L110").Value = 11950
$ws.Range("N110").Value = -20130
$ws.Range("H122").Value = 99999
$ws.Range("J122").Value = 99999
$ws.Range("L122").Value = 99999
$ws.Range("N122").Value = -109799
$ws.Range("H134").Value = 3311.625
$ws.Range("I134").Value = 3311.625
$ws.Range("K134").Value = 9934.875
$ws.Range("M134").Value = -7399.875
$ws.Range("H135").Value = 79263.63
$ws.Range("J135").Value = 79263.63
$ws.Range("L135").Value = 79263.63
$ws.Range("N135").Value = -89403.63

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4802.5
$ws.Range("I16").Value = 4287.25
$ws.Range("K16").Value = 4287.25
$ws.Range("M16").Value = -4000.25
$ws.Range("H31").Value = 3425.4736
$ws.Range("J31").Value = 1446.5
$ws.Range("L31").Value = 1446.5
$ws.Range("N31").Value = -2036.5
$ws.Range("H34").Value = 3425.4736
$ws.Range("J34").Value = 1446.5
$ws.Range("L34").Value = 1446.5
$ws.Range("N34").Value = -1850.5
$ws.Range("H107").Value = 826.0417
$ws.Range("I107").Value = 462.22223
$ws.Range("K107").Value = 462.22223
$ws.Range("M107").Value = 1457.77777
$ws.Range("H113").Value = 4802.5
$ws.Range("I113").Value = 4287.25
$ws.Range("K113").Value = 4287.25
$ws.Range("M113").Value = -2117.25
$ws.Range("H122").Value = 2112.7693
$ws.Range("I122").Value = 1959.875
$ws.Range("J122").Value = 2357.4
$ws.Range("K122").Value = 5879.625
$ws.Range("L122").Value = 7072.200000000001
$ws.Range("M122").Value = -3429.625
$ws.Range("N122").Value = -11972.2
$ws.Range("H134").Value = 1794
$ws.Range("I134").Value = 1794
$ws.Range("K134").Value = 5382
$ws.Range("M134").Value = -2847

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 16349304
$ws.Range("I2").Value = 1250108.8
$ws.Range("J2").Value = 25641118
$ws.Range("K2").Value = 7500652.800000001
$ws.Range("L2").Value = 153846708
$ws.Range("M2").Value = -7500539.800000001
$ws.Range("N2").Value = -153846934
$ws.Range("H101").Value = 12990
$ws.Range("J101").Value = 12990
$ws.Range("L101").Value = 38970
$ws.Range("N101").Value = -43838
$ws.Range("H122").Value = 1440
$ws.Range("I122").Value = 3532.6667
$ws.Range("J122").Value = 991.5714
$ws.Range("K122").Value = 31794.0003
$ws.Range("L122").Value = 8924.142600000001
$ws.Range("M122").Value = -29344.0003
$ws.Range("N122").Value = -13824.1426
$ws.Range("H138").Value = 2989.8
$ws.Range("I138").Value = 3537.25
$ws.Range("J138").Value = 800
$ws.Range("K138").Value = 10611.75
$ws.Range("L138").Value = 2400
$ws.Range("M138").Value = -5471.75
$ws.Range("N138").Value = -12680

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5146.516
$ws.Range("I80").Value = 5171.0586
$ws.Range("J80").Value = 5116.7144
$ws.Range("K80").Value = 5171.0586
$ws.Range("L80").Value = 5116.7144
$ws.Range("M80").Value = -4173.0586
$ws.Range("N80").Value = -7112.7144
$ws.Range("H83").Value = 5146.516
$ws.Range("I83").Value = 5171.0586
$ws.Range("J83").Value = 5116.7144
$ws.Range("K83").Value = 25855.293
$ws.Range("L83").Value = 25583.572
$ws.Range("M83").Value = -20863.293
$ws.Range("N83").Value = -35567.572
$ws.Range("H97").Value = 490.27585
$ws.Range("I97").Value = 345.66666
$ws.Range("K97").Value = 345.66666
$ws.Range("M97").Value = 150.33334
$ws.Range("H107").Value = 1409.9166
$ws.Range("I107").Value = 1409.9166
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1409.9166
$ws.Range("L107").Value = 0
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = 510.0834
$ws.Range("H122").Value = 4463.684
$ws.Range("I122").Value = 3988.8823
$ws.Range("J122").Value = 8499.5
$ws.Range("K122").Value = 11966.6469
$ws.Range("L122").Value = 25498.5
$ws.Range("M122").Value = -9516.6469
$ws.Range("N122").Value = -30398.5
$ws.Range("H126").Value = 4763.3213
$ws.Range("I126").Value = 3907.2307
$ws.Range("K126").Value = 11721.6921
$ws.Range("M126").Value = -9251.6921
$ws.Range("H132").Value = 6217.706
$ws.Range("I132").Value = 7375.4165
$ws.Range("J132").Value = 3439.2
$ws.Range("K132").Value = 22126.2495
$ws.Range("L132").Value = 10317.6
$ws.Range("M132").Value = -19596.2495
$ws.Range("N132").Value = -15377.6

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 171641.5
$ws.Range("I22").Value = 500250
$ws.Range("K22").Value = 500250
$ws.Range("M22").Value = -499955
$ws.Range("H27").Value = 171641.5
$ws.Range("I27").Value = 500250
$ws.Range("K27").Value = 500250
$ws.Range("M27").Value = -500143
$ws.Range("H31").Value = 1388.75
$ws.Range("I31").Value = 1388.75
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1388.75
$ws.Range("L31").Value = 0
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -1140.75
$ws.Range("H40").Value = 6314.4287
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H132").Value = 4567.8335
$ws.Range("I132").Value = 4371.7036
$ws.Range("K132").Value = 13115.1108
$ws.Range("M132").Value = -10585.1108

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2862.2856
$ws.Range("I81").Value = 1165.6923
$ws.Range("J81").Value = 5619.25
$ws.Range("K81").Value = 2331.3846
$ws.Range("L81").Value = 11238.5
$ws.Range("M81").Value = -1270.3846
$ws.Range("N81").Value = -13360.5
$ws.Range("H84").Value = 2862.2856
$ws.Range("I84").Value = 1165.6923
$ws.Range("J84").Value = 5619.25
$ws.Range("K84").Value = 11656.923
$ws.Range("L84").Value = 56192.5
$ws.Range("M84").Value = -6352.922999999999
$ws.Range("N84").Value = -66800.5
$ws.Range("H136").Value = 4622.727
$ws.Range("I136").Value = 3859.875
$ws.Range("J136").Value = 6657
$ws.Range("K136").Value = 11579.625
$ws.Range("L136").Value = 19971
$ws.Range("M136").Value = -9029.625
$ws.Range("N136").Value = -25071
